# Q3 Update - 2025
# Applies the UNHCR Mozambique (UN-MOZ) dataset refresh:
#  - The "short-url" slug (column B) changes for every data row (O60Grh -> jHo05Y)
#  - A handful of per-country statistics for Mozambique (coa) / year 2024 rows
#    (rows 430-444) are refreshed, and the "Comoros" origin-country row is
#    replaced by a new "Cote d'Ivoire" row while the remaining countries
#    shift up by one position (Ethiopia, Palestinian, Guinea).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Site-wide short-url slug update (applies to every data row, B2:B444) ---
$ws.Range("B2:B444").Value2 = "jHo05Y"

# --- Row 430 (Burundi): returned_refugees + ooc ---
$ws.Range("P430").Value2 = "12"
$ws.Range("T430").Value2 = "799"

# --- Row 433 (Dem. Rep. of the Congo): refugees / asylum_seekers / ooc ---
$ws.Range("N433").Value2 = "2936"
$ws.Range("O433").Value2 = "6404"
$ws.Range("T433").Value2 = "780"

# --- Row 434: Comoros removed, becomes Ethiopia ---
$ws.Range("F434").Value2 = "58"
$ws.Range("G434").Value2 = "Ethiopia"
$ws.Range("H434").Value2 = "ETH"
$ws.Range("I434").Value2 = "ETH"
$ws.Range("O434").Value2 = "9"

# --- Row 435: Ethiopia -> Palestinian ---
$ws.Range("F435").Value2 = "69"
$ws.Range("G435").Value2 = "Palestinian"
$ws.Range("H435").Value2 = "GAZ"
$ws.Range("I435").Value2 = "PSE"

# --- Row 436: Palestinian -> Guinea ---
$ws.Range("F436").Value2 = "79"
$ws.Range("G436").Value2 = "Guinea"
$ws.Range("H436").Value2 = "GUI"
$ws.Range("I436").Value2 = "GIN"
$ws.Range("O436").Value2 = "5"

# --- Row 437: Guinea -> Cote d'Ivoire (new entry) ---
$ws.Range("F437").Value2 = "87"
$ws.Range("G437").Value2 = "Cote d'Ivoire"
$ws.Range("H437").Value2 = "ICO"
$ws.Range("I437").Value2 = "CIV"
$ws.Range("O437").Value2 = "428"

# --- Row 439 (Rwanda): refugees / asylum_seekers / ooc ---
$ws.Range("N439").Value2 = "851"
$ws.Range("O439").Value2 = "2669"
$ws.Range("T439").Value2 = "335"

# --- Row 440 (Somalia): refugees / asylum_seekers / ooc ---
$ws.Range("N440").Value2 = "476"
$ws.Range("O440").Value2 = "1776"
$ws.Range("T440").Value2 = "90"

# --- Row 441 (Sudan): asylum_seekers / ooc ---
$ws.Range("O441").Value2 = "24"
$ws.Range("T441").Value2 = "5"

# --- Row 443 (Turkiye): asylum_seekers ---
$ws.Range("O443").Value2 = "30"

# --- Row 444 (Uganda): asylum_seekers ---
$ws.Range("O444").Value2 = "23"
